## Aggiornamento dati fino al 20/09/2021 (Formigine)
## Adds daily COVID rows 375-385 (dates 2021-09-10 .. 2021-09-20) to the
## existing time-series table on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replicate the date-column formatting (style used by A374, the last
# existing date cell) onto the new date cells A375:A385 before writing
# the values, so the new rows render exactly like the rest of column A.
$ws.Range("A374").Copy()
$ws.Range("A375:A385").PasteSpecial(-4122)

$ws.Cells.Item(375, 1).Value = 44449
$ws.Cells.Item(375, 2).Value = 3
$ws.Cells.Item(375, 3).Value = 7
$ws.Cells.Item(375, 4).Value = 20.36896933015189

$ws.Cells.Item(376, 1).Value = 44450
$ws.Cells.Item(376, 2).Value = 0
$ws.Cells.Item(376, 3).Value = 6
$ws.Cells.Item(376, 4).Value = 17.45911656870162

$ws.Cells.Item(377, 1).Value = 44451
$ws.Cells.Item(377, 2).Value = 2
$ws.Cells.Item(377, 3).Value = 7
$ws.Cells.Item(377, 4).Value = 20.36896933015189

$ws.Cells.Item(378, 1).Value = 44452
$ws.Cells.Item(378, 2).Value = 3
$ws.Cells.Item(378, 3).Value = 9
$ws.Cells.Item(378, 4).Value = 26.18867485305244

$ws.Cells.Item(379, 1).Value = 44453
$ws.Cells.Item(379, 2).Value = 1
$ws.Cells.Item(379, 3).Value = 10
$ws.Cells.Item(379, 4).Value = 29.0985276145027

$ws.Cells.Item(380, 1).Value = 44454
$ws.Cells.Item(380, 2).Value = 0
$ws.Cells.Item(380, 3).Value = 10
$ws.Cells.Item(380, 4).Value = 29.0985276145027

$ws.Cells.Item(381, 1).Value = 44455
$ws.Cells.Item(381, 2).Value = 9
$ws.Cells.Item(381, 3).Value = 18
$ws.Cells.Item(381, 4).Value = 52.37734970610487

$ws.Cells.Item(382, 1).Value = 44456
$ws.Cells.Item(382, 2).Value = 8
$ws.Cells.Item(382, 3).Value = 23
$ws.Cells.Item(382, 4).Value = 66.92661351335623

$ws.Cells.Item(383, 1).Value = 44457
$ws.Cells.Item(383, 2).Value = 4
$ws.Cells.Item(383, 3).Value = 27
$ws.Cells.Item(383, 4).Value = 78.56602455915731

$ws.Cells.Item(384, 1).Value = 44458
$ws.Cells.Item(384, 2).Value = 2
$ws.Cells.Item(384, 3).Value = 27
$ws.Cells.Item(384, 4).Value = 78.56602455915731

$ws.Cells.Item(385, 1).Value = 44459
$ws.Cells.Item(385, 2).Value = 0
$ws.Cells.Item(385, 3).Value = 24
$ws.Cells.Item(385, 4).Value = 69.8364662748065
